$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.455.20'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.16%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.653.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.25%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.54%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.57%  '

# Row 7
$ws.Range("E7").Value = '  -0.23%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.601'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.90%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.59'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.32%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.110'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.61%  '

# Row 11
$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.381'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.13%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.158'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.59%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.116.46'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.99%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.36'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.85%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '61.338.75'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.93%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000147'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.95%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.653.82'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.00%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.77'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.25%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.79'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.22%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '354.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.43%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.90'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.06%  '

# Row 22
$ws.Range("E22").Value = '  -0.10%  '

# Row 23
$ws.Range("E23").Value = '  +0.58%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.49'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.34%  '

# Row 25
$ws.Range("E25").Value = '  +3.30%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.54'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.19%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.30%  '

# Row 28
$ws.Range("E28").Value = '  +7.24%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0824'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.20%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.90'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.37%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.76'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.16%  '

# Row 32
$ws.Range("E32").Value = '  -0.09%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.53%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +13.26%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.70'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.39%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.38'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +11.09%  '

# Row 37
$ws.Range("E37").Value = '  +5.76%  '

# Row 38
$ws.Range("B38").Value = 'SuiNetwork'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.964'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +15.27%  '

# Row 39
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '339.41'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.20%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.15'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.60%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.43'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.21%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.37'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.17%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0583'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.50%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.68'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.04%  '

# Row 45
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.12'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.30%  '

# Row 46
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.631'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.80%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '135.44'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.01%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0252'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.90%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.100'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.44%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.997'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.23%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.095.75'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.66%  '
